$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "\31 52170-case-636"
$ws.Range("B3").Value = "\31 52171-case-641"
$ws.Range("C3").Value = "1 TB"
$ws.Range("D3").Value = "Black"
